$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure columns B:E are treated as plain text so values such as
# "28.208.41", "1.000", "1.0000" are preserved exactly as authored
# instead of being auto-coerced into numbers by Excel.
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range("B2").Value = 'Bitcoin'
$ws.Range("C2").Value = 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'
$ws.Range("D2").Value = '28.208.41'
$ws.Range("E2").Value = '  +1.23%  '

$ws.Range("B3").Value = 'Ethereum'
$ws.Range("C3").Value = 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'
$ws.Range("D3").Value = '1.870.13'
$ws.Range("E3").Value = '  +3.83%  '

$ws.Range("B4").Value = 'TetherUSD'
$ws.Range("C4").Value = 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -1.31%  '

$ws.Range("B5").Value = 'BNB'
$ws.Range("C5").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D5").Value = '311.62'
$ws.Range("E5").Value = '  +1.86%  '

$ws.Range("B6").Value = 'USDC'
$ws.Range("C6").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  -0.70%  '

$ws.Range("B7").Value = 'XRP'
$ws.Range("C7").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D7").Value = '0.5056'
$ws.Range("E7").Value = '  +1.18%  '

$ws.Range("B8").Value = 'Cardano'
$ws.Range("C8").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D8").Value = '0.3939'
$ws.Range("E8").Value = '  +2.45%  '

$ws.Range("B9").Value = 'Dogecoin'
$ws.Range("C9").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D9").Value = '0.09720'
$ws.Range("E9").Value = '  +3.28%  '

$ws.Range("B10").Value = 'Polygon'
$ws.Range("C10").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D10").Value = '1.143'
$ws.Range("E10").Value = '  +4.51%  '

$ws.Range("B11").Value = 'OKB'
$ws.Range("C11").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D11").Value = '40.88'
$ws.Range("E11").Value = '  +0.57%  '

$ws.Range("B12").Value = 'Polkadot'
$ws.Range("C12").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D12").Value = '6.524'
$ws.Range("E12").Value = '  +2.58%  '

$ws.Range("B13").Value = 'Solana'
$ws.Range("C13").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D13").Value = '21.00'
$ws.Range("E13").Value = '  +2.14%  '

$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.872.35'
$ws.Range("E14").Value = '  +2.89%  '

$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D15").Value = '7.435'
$ws.Range("E15").Value = '  +2.43%  '

$ws.Range("B16").Value = 'BinanceUSD'
$ws.Range("C16").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D16").Value = '1.000'
$ws.Range("E16").Value = '  -1.41%  '

$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").Value = '0.00001130'
$ws.Range("E17").Value = '  +1.91%  '

$ws.Range("B18").Value = 'Litecoin'
$ws.Range("C18").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D18").Value = '93.04'
$ws.Range("E18").Value = '  +1.17%  '

$ws.Range("B19").Value = 'TRON'
$ws.Range("C19").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D19").Value = '0.06585'
$ws.Range("E19").Value = '  +0.22%  '

$ws.Range("B20").Value = 'Avalanche'
$ws.Range("C20").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D20").Value = '17.56'
$ws.Range("E20").Value = '  +3.05%  '

$ws.Range("B21").Value = 'Dai'
$ws.Range("C21").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D21").Value = '1.0000'
$ws.Range("E21").Value = '  -0.43%  '

$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D22").Value = '6.171'
$ws.Range("E22").Value = '  +3.49%  '

$ws.Range("B23").Value = 'WrappedBTC'
$ws.Range("C23").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D23").Value = '28.275.03'
$ws.Range("E23").Value = '  +0.95%  '

$ws.Range("B24").Value = 'Cosmos'
$ws.Range("C24").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D24").Value = '11.37'
$ws.Range("E24").Value = '  +3.68%  '

$ws.Range("B25").Value = 'Toncoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D25").Value = '2.274'
$ws.Range("E25").Value = '  +1.71%  '

$ws.Range("B26").Value = 'LEO'
$ws.Range("C26").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D26").Value = '3.383'
$ws.Range("E26").Value = '  -1.33%  '

$ws.Range("B27").Value = 'LidoDAOToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D27").Value = '2.555'
$ws.Range("E27").Value = '  +8.06%  '

$ws.Range("B28").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C28").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D28").Value = '2.086.91'
$ws.Range("E28").Value = '  +2.65%  '

$ws.Range("B29").Value = 'EthereumClassic'
$ws.Range("C29").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D29").Value = '21.26'
$ws.Range("E29").Value = '  +4.15%  '

$ws.Range("B30").Value = 'Monero'
$ws.Range("C30").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D30").Value = '158.08'
$ws.Range("E30").Value = '  +0.62%  '

$ws.Range("B31").Value = 'BitcoinCash'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D31").Value = '127.71'
$ws.Range("E31").Value = '  +0.73%  '

$ws.Range("B32").Value = 'Stellar'
$ws.Range("C32").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D32").Value = '0.1060'
$ws.Range("E32").Value = '  -1.89%  '

$ws.Range("B33").Value = 'ImmutableX'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D33").Value = '1.070'
$ws.Range("E33").Value = '  +2.89%  '

$ws.Range("B34").Value = 'Filecoin'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D34").Value = '5.647'
$ws.Range("E34").Value = '  +2.32%  '

$ws.Range("B35").Value = 'HuobiToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D35").Value = '3.622'
$ws.Range("E35").Value = '  -0.28%  '

$ws.Range("B36").Value = 'FraxShare'
$ws.Range("C36").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D36").Value = '9.545'
$ws.Range("E36").Value = '  +7.23%  '

$ws.Range("B37").Value = 'Hedera'
$ws.Range("C37").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D37").Value = '0.06731'
$ws.Range("E37").Value = '  -0.94%  '

$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").Value = '0.02388'
$ws.Range("E38").Value = '  +2.95%  '

$ws.Range("B39").Value = 'Algorand'
$ws.Range("C39").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D39").Value = '0.2191'
$ws.Range("E39").Value = '  +2.69%  '

$ws.Range("B40").Value = 'Aptos'
$ws.Range("C40").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D40").Value = '11.53'
$ws.Range("E40").Value = '  +2.32%  '

$ws.Range("B41").Value = 'TheSandbox'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D41").Value = '0.6378'
$ws.Range("E41").Value = '  +4.39%  '

$ws.Range("B42").Value = 'InternetComputer(DFINITY)'
$ws.Range("C42").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D42").Value = '4.981'
$ws.Range("E42").Value = '  +0.76%  '

$ws.Range("B43").Value = 'TrustWalletToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D43").Value = '1.180'
$ws.Range("E43").Value = '  +2.67%  '

$ws.Range("B44").Value = 'Frax'
$ws.Range("C44").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D44").Value = '0.9999'
$ws.Range("E44").Value = '  -0.43%  '

$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").Value = '13.49'
$ws.Range("E45").Value = '  +2.09%  '

$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").Value = '0.6028'
$ws.Range("E46").Value = '  +2.56%  '

$ws.Range("B47").Value = 'PancakeSwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D47").Value = '3.658'
$ws.Range("E47").Value = '  -0.72%  '

$ws.Range("B48").Value = 'WEMIXTOKEN'
$ws.Range("C48").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D48").Value = '1.260'
$ws.Range("E48").Value = '  -2.52%  '

$ws.Range("B49").Value = 'NEARProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D49").Value = '1.997'
$ws.Range("E49").Value = '  +2.66%  '

$ws.Range("B50").Value = 'Quant'
$ws.Range("C50").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D50").Value = '124.10'
$ws.Range("E50").Value = '  +0.34%  '

$ws.Range("B51").Value = 'EOS'
$ws.Range("C51").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D51").Value = '1.198'
$ws.Range("E51").Value = '  +2.66%  '
